# Auto-generated: apply scheduled market-data refresh to Midgardsormr_Profits sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 946.3951
$ws.Range("I15").Value = 946.3951
$ws.Range("K15").Value = 2839.1853
$ws.Range("M15").Value = -2670.1853
$ws.Range("H69").Value = 355622
$ws.Range("I69").Value = 12977.75
$ws.Range("K69").Value = 38933.25
$ws.Range("M69").Value = -38059.25
$ws.Range("H72").Value = 355622
$ws.Range("I72").Value = 12977.75
$ws.Range("K72").Value = 116799.75
$ws.Range("M72").Value = -112431.75
$ws.Range("H88").Value = 3140.4482
$ws.Range("I88").Value = 1270.2222
$ws.Range("J88").Value = 3982.05
$ws.Range("K88").Value = 1270.2222
$ws.Range("L88").Value = 3982.05
$ws.Range("M88").Value = -864.2221999999999
$ws.Range("N88").Value = -4794.05
$ws.Range("H91").Value = 3140.4482
$ws.Range("I91").Value = 1270.2222
$ws.Range("J91").Value = 3982.05
$ws.Range("K91").Value = 1270.2222
$ws.Range("L91").Value = 3982.05
$ws.Range("M91").Value = 133.7778000000001
$ws.Range("N91").Value = -6790.05
$ws.Range("H101").Value = 1266.6666
$ws.Range("J101").Value = 1266.6666
$ws.Range("L101").Value = 3799.9998
$ws.Range("N101").Value = -7043.9998
$ws.Range("H111").Value = 2065
$ws.Range("I111").Value = 2247.5
$ws.Range("K111").Value = 6742.5
$ws.Range("M111").Value = -3675.5
$ws.Range("H125").Value = 2829
$ws.Range("I125").Value = 999
$ws.Range("J125").Value = 4659
$ws.Range("K125").Value = 8991
$ws.Range("L125").Value = 41931
$ws.Range("M125").Value = -6531
$ws.Range("N125").Value = -46851

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14657.286
$ws.Range("I32").Value = 14337.6455
$ws.Range("K32").Value = 14337.6455
$ws.Range("M32").Value = -14050.6455
$ws.Range("H61").Value = 2882.4443
$ws.Range("I61").Value = 2104.9167
$ws.Range("K61").Value = 2104.9167
$ws.Range("M61").Value = -1892.9167
$ws.Range("H88").Value = 5839.6665
$ws.Range("J88").Value = 8134.25
$ws.Range("L88").Value = 8134.25
$ws.Range("N88").Value = -8946.25
$ws.Range("H91").Value = 5839.6665
$ws.Range("J91").Value = 8134.25
$ws.Range("L91").Value = 8134.25
$ws.Range("N91").Value = -10942.25
$ws.Range("H97").Value = 2093.851
$ws.Range("I97").Value = 1715.2646
$ws.Range("J97").Value = 3084
$ws.Range("K97").Value = 1715.2646
$ws.Range("L97").Value = 3084
$ws.Range("M97").Value = -1219.2646
$ws.Range("N97").Value = -4076
$ws.Range("H110").Value = 1301
$ws.Range("I110").Value = 1301.3334
$ws.Range("K110").Value = 1301.3334
$ws.Range("M110").Value = 743.6666
$ws.Range("H122").Value = 1430.1666
$ws.Range("I122").Value = 1303.925
$ws.Range("J122").Value = 2061.375
$ws.Range("K122").Value = 3911.775
$ws.Range("L122").Value = 6184.125
$ws.Range("M122").Value = -1461.775
$ws.Range("N122").Value = -11084.125
$ws.Range("H132").Value = 2260.8
$ws.Range("I132").Value = 2260.8
$ws.Range("K132").Value = 6782.400000000001
$ws.Range("M132").Value = -4252.400000000001
$ws.Range("H136").Value = 2882.4443
$ws.Range("I136").Value = 2104.9167
$ws.Range("K136").Value = 6314.750100000001
$ws.Range("M136").Value = -3764.750100000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 515.6667
$ws.Range("I22").Value = 515.6667
$ws.Range("K22").Value = 515.6667
$ws.Range("M22").Value = -342.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 1886.4445
$ws.Range("I105").Value = 1886.4445
$ws.Range("K105").Value = 1886.4445
$ws.Range("M105").Value = -139.4445000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1228.7587
$ws.Range("I5").Value = 794.0454999999999
$ws.Range("J5").Value = 2595
$ws.Range("K5").Value = 2382.1365
$ws.Range("L5").Value = 7785
$ws.Range("M5").Value = -2270.1365
$ws.Range("N5").Value = -8009
$ws.Range("H14").Value = 249.15384
$ws.Range("I14").Value = 249.15384
$ws.Range("K14").Value = 747.4615200000001
$ws.Range("M14").Value = -574.4615200000001
$ws.Range("H86").Value = 2500
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 2500
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 7500
$ws.Range("M86").ClearContents()
$ws.Range("N86").Value = -9872
$ws.Range("H89").Value = 2500
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 2500
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 22500
$ws.Range("M89").ClearContents()
$ws.Range("N89").Value = -34356
$ws.Range("H135").Value = 1228.7587
$ws.Range("I135").Value = 794.0454999999999
$ws.Range("J135").Value = 2595
$ws.Range("K135").Value = 7146.4095
$ws.Range("L135").Value = 23355
$ws.Range("M135").Value = -4611.4095
$ws.Range("N135").Value = -28425
$ws.Range("H141").Value = 5256.9375
$ws.Range("J141").Value = 4949
$ws.Range("L141").Value = 14847
$ws.Range("N141").Value = -25207

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H19").Value = 2550251
$ws.Range("I19").Value = 67001.336
$ws.Range("K19").Value = 67001.336
$ws.Range("M19").Value = -66713.336
$ws.Range("H97").Value = 1825.1482
$ws.Range("I97").Value = 1358.5238
$ws.Range("J97").Value = 3458.3333
$ws.Range("K97").Value = 1358.5238
$ws.Range("L97").Value = 3458.3333
$ws.Range("M97").Value = -862.5237999999999
$ws.Range("N97").Value = -4450.3333
$ws.Range("H102").Value = 13104.023
$ws.Range("I102").Value = 17577.742
$ws.Range("J102").Value = 1546.9166
$ws.Range("K102").Value = 17577.742
$ws.Range("L102").Value = 1546.9166
$ws.Range("M102").Value = -15955.742
$ws.Range("N102").Value = -4790.9166
$ws.Range("H107").Value = 393.06668
$ws.Range("J107").Value = 529
$ws.Range("L107").Value = 529
$ws.Range("N107").Value = -4369
$ws.Range("H113").Value = 2613.3
$ws.Range("I113").Value = 2613.3
$ws.Range("K113").Value = 2613.3
$ws.Range("M113").Value = -443.3000000000002
$ws.Range("H122").Value = 2791.3513
$ws.Range("J122").Value = 3296.889
$ws.Range("L122").Value = 9890.667000000001
$ws.Range("N122").Value = -14790.667
$ws.Range("H126").Value = 3186
$ws.Range("I126").Value = 3498.8333
$ws.Range("J126").Value = 2810.6
$ws.Range("K126").Value = 10496.4999
$ws.Range("L126").Value = 8431.799999999999
$ws.Range("M126").Value = -8026.499899999999
$ws.Range("N126").Value = -13371.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H18").Value = 75000
$ws.Range("J18").Value = 75000
$ws.Range("L18").Value = 75000
$ws.Range("N18").Value = -75344
$ws.Range("H22").Value = 1400.3572
$ws.Range("I22").Value = 850.25
$ws.Range("J22").Value = 1620.4
$ws.Range("K22").Value = 850.25
$ws.Range("L22").Value = 1620.4
$ws.Range("M22").Value = -555.25
$ws.Range("N22").Value = -2210.4
$ws.Range("H27").Value = 1400.3572
$ws.Range("I27").Value = 850.25
$ws.Range("J27").Value = 1620.4
$ws.Range("K27").Value = 850.25
$ws.Range("L27").Value = 1620.4
$ws.Range("M27").Value = -743.25
$ws.Range("N27").Value = -1834.4
$ws.Range("H61").Value = 1221.7693
$ws.Range("I61").Value = 1327.3334
$ws.Range("J61").Value = 1131.2858
$ws.Range("K61").Value = 1327.3334
$ws.Range("L61").Value = 1131.2858
$ws.Range("M61").Value = -1125.3334
$ws.Range("N61").Value = -1535.2858
$ws.Range("H113").Value = 1221.7693
$ws.Range("I113").Value = 1327.3334
$ws.Range("J113").Value = 1131.2858
$ws.Range("K113").Value = 1327.3334
$ws.Range("L113").Value = 1131.2858
$ws.Range("M113").Value = 842.6666
$ws.Range("N113").Value = -5471.2858
$ws.Range("H132").Value = 2871.081
$ws.Range("I132").Value = 2787.7742
$ws.Range("K132").Value = 8363.3226
$ws.Range("M132").Value = -5833.3226

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H19").Value = 23451.334
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 23451.334
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 23451.334
$ws.Range("M19").ClearContents()
$ws.Range("N19").Value = -23799.334
$ws.Range("H81").Value = 3180.7273
$ws.Range("I81").Value = 3420.1
$ws.Range("J81").Value = 787
$ws.Range("K81").Value = 6840.2
$ws.Range("L81").Value = 1574
$ws.Range("M81").Value = -5779.2
$ws.Range("N81").Value = -3696
$ws.Range("H84").Value = 3180.7273
$ws.Range("I84").Value = 3420.1
$ws.Range("J84").Value = 787
$ws.Range("K84").Value = 34201
$ws.Range("L84").Value = 7870
$ws.Range("M84").Value = -28897
$ws.Range("N84").Value = -18478
